# Generate Report for Handback
# Adds a new handed-back file (f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0) as row 4
# to the Overview, zh-cn and de-de sheets, matching the existing
# "Handed back: in sync with en-US" (8a2dfb26...) entries.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
$overview.Range("B4").Value = "Handed back: in sync with en-US"
$overview.Range("C4").Value = "Handed back: in sync with en-US"

$overview.Hyperlinks.Add(
    $overview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f8b7a46d2bb7426a8d87c06c5a4a54e0handback/e2e/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
) | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: detailed status row
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
$zhcn.Range("B4").Value = ".md"
$zhcn.Range("C4").Value = "Handed back: in sync with en-US"
$zhcn.Range("D4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.zh-cn.xlf"
$zhcn.Range("E4").Value = "2016-03-22 12:02:22"
$zhcn.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("F4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
$zhcn.Range("G4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.zh-cn.xlf"
$zhcn.Range("H4").Value = "2016-03-22 12:03:00"
$zhcn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("J4").Value = "Include"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f8b7a46d2bb7426a8d87c06c5a4a54e0handback/e2e/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("D4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b44e06bbb17a0c93d068149c0342010f993760d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.zh-cn.xlf",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.zh-cn.xlf"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/f8b7a46d2bb7426a8d87c06c5a4a54e0target/e2e/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8b44e06bbb17a0c93d068149c0342010f993760d/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.zh-cn.xlf",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.zh-cn.xlf"
) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: detailed status row
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
$dede.Range("B4").Value = ".md"
$dede.Range("C4").Value = "Handed back: in sync with en-US"
$dede.Range("D4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.de-de.xlf"
$dede.Range("E4").Value = "2016-03-22 12:02:30"
$dede.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("F4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
$dede.Range("G4").Value = "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.de-de.xlf"
$dede.Range("H4").Value = "2016-03-22 12:03:16"
$dede.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("J4").Value = "Include"

$dede.Hyperlinks.Add(
    $dede.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f8b7a46d2bb7426a8d87c06c5a4a54e0handback/e2e/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("D4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8b44e06bbb17a0c93d068149c0342010f993760d/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.de-de.xlf",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.de-de.xlf"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/f8b7a46d2bb7426a8d87c06c5a4a54e0target/e2e/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8b44e06bbb17a0c93d068149c0342010f993760d/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.de-de.xlf",
    "",
    "",
    "f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0.8b44e06bbb17a0c93d068149c0342010f993760d.de-de.xlf"
) | Out-Null

Write-Output "Added handback row for f8b7a46d-2bb7-426a-8d87-c06c5a4a54e0 to Overview, zh-cn, de-de sheets"
